$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4 ---
$ws.Range("D4").Value = "Complete"
$ws.Range("E4").Value = [DateTime]"2020-07-04"
$ws.Range("F4").Value = "We forecast 10pts. Based on the scale that was used in our estimation activity.  Sprint Goal:  Display resgistered students and team names.    "

# --- Row 5 ---
$ws.Range("D5").Value = "Complete"
$ws.Range("E5").Value = [DateTime]"2020-07-04"
# F5 text unchanged

# --- Row 9 ---
$ws.Range("D9").Value = "Complete"
$ws.Range("E9").Value = [DateTime]"2020-07-04"
$ws.Range("F9").Value = "The Trello board includes both PBI's and tasks, which are linked.  The progress of all work items is represented by the list (column) they are in.  An explanation of the board layout was provided in the README."

# --- Row 10 ---
$ws.Range("D10").Value = "Complete"
$ws.Range("E10").Value = [DateTime]"2020-07-05"
$ws.Range("F10").Value = "A burndown chart was created as stated with expected, planned and actual lines.  It was updated daily."

# --- Row 11 ---
$ws.Range("D11").Value = "Complete"
$ws.Range("E11").Value = [DateTime]"2020-07-05"
$ws.Range("F11").Value = "A daily scrum has been conducted every day since 6/28/20 and documented in a log with three sections as specified."

# --- Row 12 ---
$ws.Range("D12").Value = "Complete"
$ws.Range("E12").Value = [DateTime]"2020-07-05"
$ws.Range("F12").Value = "A daily scrum has been conducted every day since 6/28/20 and documented in a log with three sections as specified."

# --- Row 13 ---
$ws.Range("D13").Value = "Complete"
$ws.Range("E13").Value = [DateTime]"2020-07-05"
$ws.Range("F13").Value = "A daily scrum has been conducted every day since 6/28/20 and documented in a log with three sections as specified."

# --- Row 14 ---
$ws.Range("D14").Value = "Complete"
$ws.Range("E14").Value = [DateTime]"2020-07-05"
$ws.Range("F14").Value = "A daily scrum has been conducted every day since 6/28/20 and documented in a log with three sections as specified."

# --- Row 15 ---
$ws.Range("D15").Value = "Complete"
$ws.Range("E15").Value = [DateTime]"2020-07-04"
$ws.Range("F15").Value = "The burndown chart is in the repo.  URLs for the burndown and the task board are in the README"

# --- Row 16 ---
$ws.Range("D16").Value = "Complete"
$ws.Range("E16").Value = [DateTime]"2020-07-05"
$ws.Range("F16").Value = "We mob programmed for the entire Sprint. Several photos are in the repo.  Links and photos are in the README."

# --- Row 17 ---
$ws.Range("D17").Value = "Complete"
$ws.Range("E17").Value = [DateTime]"2020-07-05"
$ws.Range("F17").Value = "We used TDD for the entire sprint.  We have atotal of 21 tests written with 52 assertions that all pass.  Anu confirmed the evidence in our README met requirements."

# --- Row 18 --- (D18 stays "Pending")
$ws.Range("F18").Value = "Scheduled for 7/6/20, 1PM"

# --- Row 19 ---
$ws.Range("D19").Value = "Complete"
$ws.Range("E19").Value = [DateTime]"2020-07-05"
$ws.Range("F19").Value = "Our app was pushed to the production server and works well."
$ws.Rows.Item(19).RowHeight = 28.8

# --- Row 20 --- (D20 stays "Pending")
$ws.Range("F20").Value = "Richard is scheduled to meet with us at our Sprint Review on'7/6/20, 1PM"

# --- Selection change (bottomRight pane now at F1) ---
$ws.Range("F1").Select()
